$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need an explicit text format
# so Excel stores them as text (preserving exact formatting/precision)
# instead of silently converting them to numbers.
$textCells = @("D5", "D6", "D24", "D27", "D28", "D30", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D42", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '76.318.12'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '2.864.12'
$ws.Range('E3').Value = '  +7.85%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '194.60'
$ws.Range('E5').Value = '  +4.78%  '
$ws.Range('D6').Value = '600.46'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +3.77%  '
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('D10').Value = '2.861.60'
$ws.Range('E10').Value = '  +7.80%  '
$ws.Range('E11').Value = '  +10.67%  '
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').Value = '3.388.31'
$ws.Range('E14').Value = '  +7.78%  '
$ws.Range('D15').Value = '76.034.54'
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('E16').Value = '  +4.46%  '
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').Value = '2.863.01'
$ws.Range('E18').Value = '  +7.96%  '
$ws.Range('E19').Value = '  -1.84%  '
$ws.Range('E20').Value = '  +5.12%  '
$ws.Range('E21').Value = '  +3.19%  '
$ws.Range('E22').Value = '  +4.32%  '
$ws.Range('E23').Value = '  +2.17%  '
$ws.Range('D24').Value = '72.01'
$ws.Range('E24').Value = '  +4.06%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '3.010.79'
$ws.Range('E26').Value = '  +7.74%  '
$ws.Range('D27').Value = '4.23'
$ws.Range('E27').Value = '  +2.61%  '
$ws.Range('D28').Value = '9.76'
$ws.Range('E28').Value = '  +4.82%  '
$ws.Range('E29').Value = '  +12.76%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').Value = '516.47'
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').Value = '7.73'
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('E34').Value = '  +4.98%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '167.02'
$ws.Range('E36').Value = '  +2.19%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '20.06'
$ws.Range('E38').Value = '  +4.71%  '
$ws.Range('D39').Value = '19.47'
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').Value = '186.13'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Value = '0.345'
$ws.Range('E42').Value = '  +5.84%  '
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('D45').Value = '1.24'
$ws.Range('E45').Value = '  +5.22%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0897'
$ws.Range('E46').Value = '  +6.58%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '40.33'
$ws.Range('E47').Value = '  +3.50%  '
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('E49').Value = '  +10.51%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.668'
$ws.Range('E50').Value = '  +13.47%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '3.76'
$ws.Range('E51').Value = '  +4.14%  '
